# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) across the leve-profit sheets with newly pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4185.684
$ws.Range("J40").Value = 5890.25
$ws.Range("L40").Value = 5890.25
$ws.Range("N40").Value = -6240.25
$ws.Range("H64").Value = 3812.1428
$ws.Range("I64").Value = 3112.1924
$ws.Range("J64").Value = 5834.222
$ws.Range("K64").Value = 3112.1924
$ws.Range("L64").Value = 5834.222
$ws.Range("M64").Value = -2864.1924
$ws.Range("N64").Value = -6330.222
$ws.Range("H67").Value = 3812.1428
$ws.Range("I67").Value = 3112.1924
$ws.Range("J67").Value = 5834.222
$ws.Range("K67").Value = 3112.1924
$ws.Range("L67").Value = 5834.222
$ws.Range("M67").Value = -2254.1924
$ws.Range("N67").Value = -7550.222
$ws.Range("H69").Value = 2870
$ws.Range("I69").Value = 2600
$ws.Range("K69").Value = 7800
$ws.Range("M69").Value = -6926
$ws.Range("H72").Value = 2870
$ws.Range("I72").Value = 2600
$ws.Range("K72").Value = 23400
$ws.Range("M72").Value = -19032
$ws.Range("H98").Value = 607.64105
$ws.Range("I98").Value = 453.2069
$ws.Range("J98").Value = 1055.5
$ws.Range("K98").Value = 453.2069
$ws.Range("L98").Value = 1055.5
$ws.Range("M98").Value = 1044.7931
$ws.Range("N98").Value = -4051.5
$ws.Range("H103").Value = 655.7143
$ws.Range("I103").Value = 475
$ws.Range("J103").Value = 728
$ws.Range("K103").Value = 1425
$ws.Range("L103").Value = 2184
$ws.Range("M103").Value = -839
$ws.Range("N103").Value = -3356
$ws.Range("H112").Value = 1573.5294
$ws.Range("I112").Value = 1063.3334
$ws.Range("J112").Value = 1682.8572
$ws.Range("K112").Value = 3190.0002
$ws.Range("L112").Value = 5048.571599999999
$ws.Range("M112").Value = -2082.0002
$ws.Range("N112").Value = -7264.571599999999
$ws.Range("H113").Value = 3876.5
$ws.Range("I113").Value = 2920
$ws.Range("J113").Value = 4391.5386
$ws.Range("K113").Value = 2920
$ws.Range("L113").Value = 4391.5386
$ws.Range("M113").Value = 334
$ws.Range("N113").Value = -10899.5386
$ws.Range("H122").Value = 607.64105
$ws.Range("I122").Value = 453.2069
$ws.Range("J122").Value = 1055.5
$ws.Range("K122").Value = 1359.6207
$ws.Range("L122").Value = 3166.5
$ws.Range("M122").Value = 1090.3793
$ws.Range("N122").Value = -8066.5
$ws.Range("H138").Value = 1732.7971
$ws.Range("I138").Value = 1272.4412
$ws.Range("J138").Value = 2180
$ws.Range("K138").Value = 3817.3236
$ws.Range("L138").Value = 6540
$ws.Range("M138").Value = 1322.6764
$ws.Range("N138").Value = -16820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1764.6111
$ws.Range("I2").Value = 1947.3572
$ws.Range("J2").Value = 1125
$ws.Range("K2").Value = 1947.3572
$ws.Range("L2").Value = 1125
$ws.Range("M2").Value = -1834.3572
$ws.Range("N2").Value = -1351
$ws.Range("H45").Value = 1357
$ws.Range("I45").Value = 1148
$ws.Range("K45").Value = 1148
$ws.Range("M45").Value = -771
$ws.Range("H63").Value = 4983.4443
$ws.Range("I63").Value = 10900
$ws.Range("J63").Value = 2025.1666
$ws.Range("K63").Value = 10900
$ws.Range("L63").Value = 2025.1666
$ws.Range("M63").Value = -10214
$ws.Range("N63").Value = -3397.1666
$ws.Range("H66").Value = 4983.4443
$ws.Range("I66").Value = 10900
$ws.Range("J66").Value = 2025.1666
$ws.Range("K66").Value = 54500
$ws.Range("L66").Value = 10125.833
$ws.Range("M66").Value = -51068
$ws.Range("N66").Value = -16989.833
$ws.Range("H74").Value = 2051
$ws.Range("I74").Value = 1070.3793
$ws.Range("J74").Value = 4082.2856
$ws.Range("K74").Value = 1070.3793
$ws.Range("L74").Value = 4082.2856
$ws.Range("M74").Value = -196.3793000000001
$ws.Range("N74").Value = -5830.2856
$ws.Range("H77").Value = 2051
$ws.Range("I77").Value = 1070.3793
$ws.Range("J77").Value = 4082.2856
$ws.Range("K77").Value = 5351.896500000001
$ws.Range("L77").Value = 20411.428
$ws.Range("M77").Value = -983.8965000000007
$ws.Range("N77").Value = -29147.428
$ws.Range("H116").Value = 1764.6111
$ws.Range("I116").Value = 1947.3572
$ws.Range("J116").Value = 1125
$ws.Range("K116").Value = 1947.3572
$ws.Range("L116").Value = 1125
$ws.Range("M116").Value = 346.6428000000001
$ws.Range("N116").Value = -5713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1764.6111
$ws.Range("I3").Value = 1947.3572
$ws.Range("J3").Value = 1125
$ws.Range("K3").Value = 1947.3572
$ws.Range("L3").Value = 1125
$ws.Range("M3").Value = -1833.3572
$ws.Range("N3").Value = -1353
$ws.Range("H20").Value = 1834.8334
$ws.Range("I20").Value = 1808.9231
$ws.Range("J20").Value = 1865.4546
$ws.Range("K20").Value = 1808.9231
$ws.Range("L20").Value = 1865.4546
$ws.Range("M20").Value = -1561.9231
$ws.Range("N20").Value = -2359.4546
$ws.Range("H75").Value = 13902.333
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21872
$ws.Range("H78").Value = 13902.333
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69360
$ws.Range("H86").Value = 6341.864
$ws.Range("I86").Value = 5825.9165
$ws.Range("J86").Value = 6961
$ws.Range("K86").Value = 5825.9165
$ws.Range("L86").Value = 6961
$ws.Range("M86").Value = -4702.9165
$ws.Range("N86").Value = -9207
$ws.Range("H89").Value = 6341.864
$ws.Range("I89").Value = 5825.9165
$ws.Range("J89").Value = 6961
$ws.Range("K89").Value = 29129.5825
$ws.Range("L89").Value = 34805
$ws.Range("M89").Value = -23513.5825
$ws.Range("N89").Value = -46037
$ws.Range("H94").Value = 1002.1818
$ws.Range("I94").Value = 773.0769
$ws.Range("J94").Value = 1333.1111
$ws.Range("K94").Value = 773.0769
$ws.Range("L94").Value = 1333.1111
$ws.Range("M94").Value = -322.0769
$ws.Range("N94").Value = -2235.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2390.8572
$ws.Range("I31").Value = 1419.6923
$ws.Range("J31").Value = 3969
$ws.Range("K31").Value = 1419.6923
$ws.Range("L31").Value = 3969
$ws.Range("M31").Value = -1124.6923
$ws.Range("N31").Value = -4559
$ws.Range("H34").Value = 2390.8572
$ws.Range("I34").Value = 1419.6923
$ws.Range("J34").Value = 3969
$ws.Range("K34").Value = 1419.6923
$ws.Range("L34").Value = 3969
$ws.Range("M34").Value = -1217.6923
$ws.Range("N34").Value = -4373
$ws.Range("H62").Value = 3089506
$ws.Range("I62").Value = 4631893
$ws.Range("J62").Value = 4732.3335
$ws.Range("K62").Value = 4631893
$ws.Range("L62").Value = 4732.3335
$ws.Range("M62").Value = -4631269
$ws.Range("N62").Value = -5980.3335
$ws.Range("H65").Value = 3089506
$ws.Range("I65").Value = 4631893
$ws.Range("J65").Value = 4732.3335
$ws.Range("K65").Value = 23159465
$ws.Range("L65").Value = 23661.6675
$ws.Range("M65").Value = -23156345
$ws.Range("N65").Value = -29901.6675
$ws.Range("H132").Value = 2450.8147
$ws.Range("I132").Value = 1398
$ws.Range("J132").Value = 3070.1177
$ws.Range("K132").Value = 4194
$ws.Range("L132").Value = 9210.3531
$ws.Range("M132").Value = -1664
$ws.Range("N132").Value = -14270.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 2940
$ws.Range("I56").Value = 2940
$ws.Range("K56").Value = 2940
$ws.Range("M56").Value = -2410
$ws.Range("H97").Value = 746
$ws.Range("I97").Value = 119
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 357
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = 139
$ws.Range("N97").Value = -6992
$ws.Range("H116").Value = 4299.6113
$ws.Range("I116").Value = 770.4286
$ws.Range("J116").Value = 6545.4546
$ws.Range("K116").Value = 2311.2858
$ws.Range("L116").Value = 19636.3638
$ws.Range("M116").Value = 1130.7142
$ws.Range("N116").Value = -26520.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 71433890
$ws.Range("I100").Value = 9849.833000000001
$ws.Range("K100").Value = 9849.833000000001
$ws.Range("M100").Value = -9308.833000000001
$ws.Range("H122").Value = 3182.7693
$ws.Range("I122").Value = 2875
$ws.Range("J122").Value = 3319.5557
$ws.Range("K122").Value = 8625
$ws.Range("L122").Value = 9958.667099999999
$ws.Range("M122").Value = -6175
$ws.Range("N122").Value = -14858.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 19400
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 24200
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 24200
$ws.Range("M82").Value = -4617
$ws.Range("N82").Value = -24966
$ws.Range("H85").Value = 19400
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 24200
$ws.Range("K85").Value = 5000
$ws.Range("L85").Value = 24200
$ws.Range("M85").Value = -3674
$ws.Range("N85").Value = -26852
$ws.Range("H132").Value = 1703.4857
$ws.Range("I132").Value = 762.8611
$ws.Range("J132").Value = 2699.4412
$ws.Range("K132").Value = 2288.5833
$ws.Range("L132").Value = 8098.323600000001
$ws.Range("M132").Value = 241.4167000000002
$ws.Range("N132").Value = -13158.3236
$ws.Range("H136").Value = 23393670
$ws.Range("I136").Value = 28572754
$ws.Range("J136").Value = 15154220
$ws.Range("K136").Value = 85718262
$ws.Range("L136").Value = 45462660
$ws.Range("M136").Value = -85715712
$ws.Range("N136").Value = -45467760
